$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 2 (Louise Deforge / saidunuhu579@gmail.com / ...),
# shifting row 3 (Janice Holley / ...) up into row 2.
$ws.Rows.Item(2).Delete()

# The row that is now row 2 (previously row 3, Janice Holley) gets a freshly
# randomized proxy assignment.
$ws.Range("E2").Value = "198.140.141.18:47299"
$ws.Range("F2").Value = "xL50iQ642EZkOn"
$ws.Range("G2").Value = "i8IhO2946aEDpf"

# Update the selection to match the new state of the sheet -- the user's
# selection still spans down to row 3 even though the data now ends at row 2.
$ws.Range("A2:G3").Select() | Out-Null
